$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (price & volume/1h changes, plus two row re-orderings
# caused by ranking swaps) as captured by the GitHub Actions refresh commit.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.773.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.275.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.00"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.75"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.425"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0950"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.14"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.612.43"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.75"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.84"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +9.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.80"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.815"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.268.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.656.44"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0937"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.00%  "
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.25"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.46%  "
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.99"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.57"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.36"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.87"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "170.52"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.140"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.60"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.20%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.66"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.10"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.35%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0660"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.23%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.41"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.66%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.47"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.77%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0250"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.18%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.77"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000227"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -10.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.100"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.03%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.49"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.72%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.22"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.27"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.473.19"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.65"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.27"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +8.37%  "
